# Remove the first data row (2025-10-09) from the "Chart" sheet.
# This shifts all subsequent rows up by one, which matches the
# target edit: the oldest date row is dropped and every later row's
# Not indexed / Indexed / Impressions figures move up one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows("2:2").Delete()
